$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 8
$ws.Range("H8").Value = 175
$ws.Range("I8").Value = 175
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 525
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -386
$ws.Range("N8").ClearContents()

# row 127
$ws.Range("H127").Value = 924.35297
$ws.Range("I127").Value = 528.75
$ws.Range("J127").Value = 1276
$ws.Range("K127").Value = 1586.25
$ws.Range("L127").Value = 3828
$ws.Range("M127").Value = 3373.75
$ws.Range("N127").Value = -13748

# row 138
$ws.Range("H138").Value = 2594.94
$ws.Range("I138").Value = 2807.1428
$ws.Range("J138").Value = 2512.4167
$ws.Range("K138").Value = 8421.428400000001
$ws.Range("L138").Value = 7537.250100000001
$ws.Range("M138").Value = -3281.428400000001
$ws.Range("N138").Value = -17817.2501


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 19
$ws.Range("H19").Value = 1900
$ws.Range("I19").Value = 1900
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1900
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -1671
$ws.Range("N19").ClearContents()

# row 61
$ws.Range("H61").Value = 9806133
$ws.Range("I61").Value = 20834674
$ws.Range("J61").Value = 2984.5557
$ws.Range("K61").Value = 20834674
$ws.Range("L61").Value = 2984.5557
$ws.Range("M61").Value = -20834462
$ws.Range("N61").Value = -3408.5557

# row 132
$ws.Range("H132").Value = 5315.3105
$ws.Range("I132").Value = 6255.0527
$ws.Range("J132").Value = 3529.8
$ws.Range("K132").Value = 18765.1581
$ws.Range("L132").Value = 10589.4
$ws.Range("M132").Value = -16235.1581
$ws.Range("N132").Value = -15649.4

# row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# row 136
$ws.Range("H136").Value = 9806133
$ws.Range("I136").Value = 20834674
$ws.Range("J136").Value = 2984.5557
$ws.Range("K136").Value = 62504022
$ws.Range("L136").Value = 8953.667099999999
$ws.Range("M136").Value = -62501472
$ws.Range("N136").Value = -14053.6671


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# row 134
$ws.Range("H134").Value = 4235.643
$ws.Range("I134").Value = 4883.1665
$ws.Range("J134").Value = 3750
$ws.Range("K134").Value = 14649.4995
$ws.Range("L134").Value = 11250
$ws.Range("M134").Value = -12114.4995
$ws.Range("N134").Value = -16320

# row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 15
$ws.Range("H15").Value = 14999.5
$ws.Range("I15").Value = 5000
$ws.Range("J15").Value = 19999.25
$ws.Range("K15").Value = 5000
$ws.Range("L15").Value = 19999.25
$ws.Range("M15").Value = -4830
$ws.Range("N15").Value = -20339.25

# row 31
$ws.Range("H31").Value = 3481.2407
$ws.Range("I31").Value = 1585.1333
$ws.Range("J31").Value = 4210.5127
$ws.Range("K31").Value = 1585.1333
$ws.Range("L31").Value = 4210.5127
$ws.Range("M31").Value = -1290.1333
$ws.Range("N31").Value = -4800.5127

# row 34
$ws.Range("H34").Value = 3481.2407
$ws.Range("I34").Value = 1585.1333
$ws.Range("J34").Value = 4210.5127
$ws.Range("K34").Value = 1585.1333
$ws.Range("L34").Value = 4210.5127
$ws.Range("M34").Value = -1383.1333
$ws.Range("N34").Value = -4614.5127

# row 51
$ws.Range("H51").Value = 18349.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 18349.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 18349.5
$ws.Range("N51").Value = -19821.5

# row 59
$ws.Range("H59").Value = 23982.143
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 23982.143
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 23982.143
$ws.Range("N59").Value = -26272.143

# row 60
$ws.Range("H60").Value = 10747
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 10747
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 10747
$ws.Range("N60").Value = -11769

# row 61
$ws.Range("H61").Value = 18349.5
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 18349.5
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 18349.5
$ws.Range("N61").Value = -19045.5

# row 68
$ws.Range("H68").Value = 28199.4
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 28199.4
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 28199.4
$ws.Range("N68").Value = -29697.4

# row 71
$ws.Range("H71").Value = 28199.4
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 28199.4
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 84598.20000000001
$ws.Range("N71").Value = -92086.20000000001

# row 74
$ws.Range("H74").Value = 28249.75
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 28249.75
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 28249.75
$ws.Range("N74").Value = -29997.75

# row 77
$ws.Range("H77").Value = 28249.75
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 28249.75
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 84749.25
$ws.Range("N77").Value = -93485.25


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 60
$ws.Range("H60").Value = 2556.561
$ws.Range("I60").Value = 423.16666
$ws.Range("J60").Value = 2922.2856
$ws.Range("K60").Value = 1269.49998
$ws.Range("L60").Value = 8766.856800000001
$ws.Range("M60").Value = -1018.49998
$ws.Range("N60").Value = -9268.856800000001

# row 68
$ws.Range("H68").Value = 1394.7084
$ws.Range("I68").Value = 769.4138
$ws.Range("J68").Value = 1665.3582
$ws.Range("K68").Value = 2308.2414
$ws.Range("L68").Value = 4996.0746
$ws.Range("M68").Value = -1497.2414
$ws.Range("N68").Value = -6618.0746

# row 71
$ws.Range("H71").Value = 1394.7084
$ws.Range("I71").Value = 769.4138
$ws.Range("J71").Value = 1665.3582
$ws.Range("K71").Value = 6924.724200000001
$ws.Range("L71").Value = 14988.2238
$ws.Range("M71").Value = -2868.724200000001
$ws.Range("N71").Value = -23100.2238

# row 107
$ws.Range("H107").Value = 1575.1594
$ws.Range("I107").Value = 314.4516
$ws.Range("J107").Value = 2603.6316
$ws.Range("K107").Value = 943.3548
$ws.Range("L107").Value = 7810.8948
$ws.Range("M107").Value = 976.6452
$ws.Range("N107").Value = -11650.8948

# row 122
$ws.Range("H122").Value = 5782.737
$ws.Range("I122").Value = 382.76923
$ws.Range("J122").Value = 17482.666
$ws.Range("K122").Value = 3444.92307
$ws.Range("L122").Value = 157343.994
$ws.Range("M122").Value = -994.9230699999998

# row 131
$ws.Range("H131").Value = 1184.5454
$ws.Range("I131").Value = 1030
$ws.Range("J131").Value = 1200
$ws.Range("K131").Value = 3090
$ws.Range("L131").Value = 3600
$ws.Range("M131").Value = 1950
$ws.Range("N131").Value = -13680


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 3
$ws.Range("H3").Value = 2666858.5
$ws.Range("I3").Value = 3000230.2
$ws.Range("J3").Value = 1000000
$ws.Range("K3").Value = 3000230.2
$ws.Range("L3").Value = 1000000
$ws.Range("M3").Value = -3000114.2
$ws.Range("N3").Value = -1000232

# row 132
$ws.Range("H132").Value = 2662.8147
$ws.Range("I132").Value = 2146.2856
$ws.Range("J132").Value = 3219.077
$ws.Range("K132").Value = 6438.8568
$ws.Range("L132").Value = 9657.231
$ws.Range("M132").Value = -3908.8568
$ws.Range("N132").Value = -14717.231


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 7002.7646
$ws.Range("I22").Value = 928.5714
$ws.Range("J22").Value = 11254.7
$ws.Range("K22").Value = 928.5714
$ws.Range("L22").Value = 11254.7
$ws.Range("M22").Value = -633.5714
$ws.Range("N22").Value = -11844.7

# row 27
$ws.Range("H27").Value = 7002.7646
$ws.Range("I27").Value = 928.5714
$ws.Range("J27").Value = 11254.7
$ws.Range("K27").Value = 928.5714
$ws.Range("L27").Value = 11254.7
$ws.Range("M27").Value = -821.5714
$ws.Range("N27").Value = -11468.7

# row 36
$ws.Range("H36").Value = 100715
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 100715
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 100715
$ws.Range("N36").Value = -101839

# row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# row 136
$ws.Range("H136").Value = 3969253.5
$ws.Range("I136").Value = 1154.8636
$ws.Range("J136").Value = 8334162
$ws.Range("K136").Value = 3464.5908
$ws.Range("L136").Value = 25002486
$ws.Range("M136").Value = -914.5907999999999
$ws.Range("N136").Value = -25007586

# row 137
$ws.Range("H137").Value = 120000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 120000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 120000
$ws.Range("N137").Value = -130200


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 136
$ws.Range("H136").Value = 2959.1592
$ws.Range("I136").Value = 2609.0588
$ws.Range("J136").Value = 4149.5
$ws.Range("K136").Value = 7827.176399999999
$ws.Range("L136").Value = 12448.5
$ws.Range("M136").Value = -5277.176399999999
$ws.Range("N136").Value = -17548.5

